# Update "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamps for the e39f574c-... row after
# regenerating the handback report.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: Latest HO Xliff Generate Date
$wsOverview.Range("G4").Value = "2016-09-07 03:10:56"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime
$wsZhCn.Range("H4").Value = "2016-09-07 03:10:44"
$wsZhCn.Range("K4").Value = "2016-09-07 03:11:56"

# de-de sheet: Correspond Handoff Datetime / Correspond Handback DateTime
$wsDeDe.Range("H4").Value = "2016-09-07 03:10:56"
$wsDeDe.Range("K4").Value = "2016-09-07 03:12:19"
